# iron-skillet "build all" refresh on the "set commands" sheet:
#  - bump the iron-skillet-version tag comment from 0.0.1 to 0.0.2
#  - add a new "Executable Linked Format" mlav-engine-filebased-enabled
#    set-command line to each of the 5 antivirus profiles
#    (Alert-Only-AV, Outbound-AV, Inbound-AV, Internal-AV, Exception-AV),
#    right after each profile's existing "PowerShell Script 2" line.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("set commands")

function Insert-AfterMatch($precedingRow, $mustContain, $newText) {
    # $precedingRow is the row that must already contain $mustContain
    # (the last "PowerShell Script 2" line of a profile). The new row is
    # inserted immediately below it, pushing that row and everything
    # after it down by one.
    $precedingCell = $ws.Cells.Item($precedingRow, 1)
    $current = $precedingCell.Text
    if ($current -notlike "*$mustContain*") {
        throw "Expected row $precedingRow to contain '$mustContain' but found: $current"
    }
    $newRow = $precedingRow + 1
    $target = $ws.Cells.Item($newRow, 1)
    $target.EntireRow.Insert()
    $target.Value = $newText
}

# 1. Bump the version comment string in place (same row, no shift).
$verCell = $ws.Range("A167")
$verText = $verCell.Text
if ($verText -notlike "*iron-skillet-version*0.0.1*") {
    throw "Expected A167 to hold the 0.0.1 iron-skillet-version tag but found: $verText"
}
$verCell.Value = 'set tag iron-skillet-version comments ""version 0.0.2 for 10.1: version of this IronSkillet template file""'

# 2. Insert the five new "Executable Linked Format" rows, working from the
#    bottom-most insertion point upward so that the row numbers below
#    (taken from the original layout) are not invalidated by earlier inserts.

Insert-AfterMatch 537 "Exception-AV mlav-engine-filebased-enabled" `
    'set profiles virus Exception-AV mlav-engine-filebased-enabled ""Executable Linked Format"" mlav-policy-action enable'

Insert-AfterMatch 513 "Internal-AV mlav-engine-filebased-enabled" `
    'set profiles virus Internal-AV mlav-engine-filebased-enabled ""Executable Linked Format"" mlav-policy-action enable'

Insert-AfterMatch 489 "Inbound-AV mlav-engine-filebased-enabled" `
    'set profiles virus Inbound-AV mlav-engine-filebased-enabled ""Executable Linked Format"" mlav-policy-action enable'

Insert-AfterMatch 465 "Outbound-AV mlav-engine-filebased-enabled" `
    'set profiles virus Outbound-AV mlav-engine-filebased-enabled ""Executable Linked Format"" mlav-policy-action enable'

Insert-AfterMatch 441 "Alert-Only-AV mlav-engine-filebased-enabled" `
    'set profiles virus Alert-Only-AV mlav-engine-filebased-enabled ""Executable Linked Format"" mlav-policy-action enable(alert-only)'
